# Update crypto price/volume columns per the latest data refresh.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the original inline-string cells) rather than re-parsing
# strings such as "1.014" or "27.565.52" as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.565.52"
$ws.Range("E2").Value = "'  +1.98%  "
$ws.Range("D3").Value = "'1.865.02"
$ws.Range("E3").Value = "'  +0.96%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("D5").Value = "'312.48"
$ws.Range("E5").Value = "'  +0.89%  "
$ws.Range("E6").Value = "'  -0.20%  "
$ws.Range("D7").Value = "'0.4783"
$ws.Range("E7").Value = "'  +0.37%  "
$ws.Range("E8").Value = "'  +3.35%  "
$ws.Range("D9").Value = "'0.07348"
$ws.Range("E9").Value = "'  +1.55%  "
$ws.Range("D10").Value = "'0.9345"
$ws.Range("E10").Value = "'  +0.31%  "
$ws.Range("D11").Value = "'20.81"
$ws.Range("E11").Value = "'  +4.72%  "
$ws.Range("D12").Value = "'0.07799"
$ws.Range("D13").Value = "'1.897.48"
$ws.Range("E13").Value = "'  +2.41%  "
$ws.Range("D14").Value = "'5.444"
$ws.Range("D15").Value = "'6.570"
$ws.Range("D16").Value = "'90.37"
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = "'  -0.24%  "
$ws.Range("D18").Value = "'0.000008831"
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("D20").Value = "'27.605.97"
$ws.Range("E20").Value = "'  +2.00%  "
$ws.Range("D21").Value = "'14.67"
$ws.Range("E21").Value = "'  +0.85%  "
$ws.Range("D22").Value = "'5.103"
$ws.Range("E22").Value = "'  +0.99%  "
$ws.Range("E23").Value = "'  +0.67%  "
$ws.Range("D24").Value = "'1.937"
$ws.Range("E24").Value = "'  +0.67%  "
$ws.Range("D25").Value = "'156.22"
$ws.Range("E25").Value = "'  +2.23%  "
$ws.Range("E26").Value = "'  +1.14%  "
$ws.Range("D27").Value = "'2.028"
$ws.Range("E27").Value = "'  +1.76%  "
$ws.Range("D28").Value = "'115.52"
$ws.Range("E28").Value = "'  +0.78%  "
$ws.Range("D29").Value = "'4.949"
$ws.Range("E29").Value = "'  +0.14%  "
$ws.Range("D30").Value = "'0.08885"
$ws.Range("E30").Value = "'  +0.19%  "
$ws.Range("E31").Value = "'  +0.36%  "
$ws.Range("E32").Value = "'  +2.59%  "
$ws.Range("D33").Value = "'0.7594"
$ws.Range("E33").Value = "'  +2.86%  "
$ws.Range("D34").Value = "'4.601"
$ws.Range("E34").Value = "'  +1.97%  "
$ws.Range("D35").Value = "'2.691"
$ws.Range("E35").Value = "'  +1.11%  "
$ws.Range("D36").Value = "'1.132"
$ws.Range("E36").Value = "'  +1.57%  "
$ws.Range("E37").Value = "'  +3.19%  "
$ws.Range("D38").Value = "'0.5677"
$ws.Range("E38").Value = "'  +7.46%  "
$ws.Range("D39").Value = "'0.05359"
$ws.Range("E39").Value = "'  +1.91%  "
$ws.Range("D40").Value = "'2.981"
$ws.Range("E40").Value = "'  +0.46%  "
$ws.Range("D41").Value = "'7.042"
$ws.Range("E41").Value = "'  +0.13%  "
$ws.Range("D42").Value = "'8.532"
$ws.Range("E42").Value = "'  +2.98%  "
$ws.Range("D43").Value = "'0.1527"
$ws.Range("E43").Value = "'  +0.52%  "
$ws.Range("D44").Value = "'0.4889"
$ws.Range("E44").Value = "'  +3.04%  "
$ws.Range("D45").Value = "'10.67"
$ws.Range("E45").Value = "'  +0.50%  "
$ws.Range("D46").Value = "'105.39"
$ws.Range("E46").Value = "'  +3.45%  "
$ws.Range("E47").Value = "'  -0.25%  "
$ws.Range("D48").Value = "'1.664"
$ws.Range("E48").Value = "'  +3.15%  "
$ws.Range("E49").Value = "'  +2.52%  "
$ws.Range("E50").Value = "'  +0.55%  "
$ws.Range("D51").Value = "'0.9114"
$ws.Range("E51").Value = "'  +2.09%  "
